# Auto-generated edit script: refresh market-price derived columns (H:N)
# for the scheduled-runner update of Sheets (per commit "chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4433.88
$ws.Range("I19").Value = 5618.1577
$ws.Range("J19").Value = 683.6667
$ws.Range("K19").Value = 5618.1577
$ws.Range("L19").Value = 683.6667
$ws.Range("M19").Value = -5443.1577
$ws.Range("N19").Value = -1033.6667

$ws.Range("H112").Value = 1711.2174
$ws.Range("I112").Value = 892
$ws.Range("J112").Value = 1938.7778
$ws.Range("K112").Value = 2676
$ws.Range("L112").Value = 5816.3334
$ws.Range("M112").Value = -1568
$ws.Range("N112").Value = -8032.3334

$ws.Range("H116").Value = 59864.668
$ws.Range("J116").Value = 3849.75
$ws.Range("L116").Value = 3849.75
$ws.Range("N116").Value = -10733.75

$ws.Range("H138").Value = 1454.6522
$ws.Range("I138").Value = 540.4375
$ws.Range("J138").Value = 2245.3242
$ws.Range("K138").Value = 1621.3125
$ws.Range("L138").Value = 6735.9726
$ws.Range("M138").Value = 3518.6875
$ws.Range("N138").Value = -17015.9726

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1090.7273
$ws.Range("I45").Value = 897.7143
$ws.Range("J45").Value = 1428.5
$ws.Range("K45").Value = 897.7143
$ws.Range("L45").Value = 1428.5
$ws.Range("M45").Value = -520.7143
$ws.Range("N45").Value = -2182.5

$ws.Range("H88").Value = 2189.2307
$ws.Range("I88").Value = 2307.2
$ws.Range("J88").Value = 1978.5714
$ws.Range("K88").Value = 2307.2
$ws.Range("L88").Value = 1978.5714
$ws.Range("M88").Value = -1901.2
$ws.Range("N88").Value = -2790.5714

$ws.Range("H91").Value = 2189.2307
$ws.Range("I91").Value = 2307.2
$ws.Range("J91").Value = 1978.5714
$ws.Range("K91").Value = 2307.2
$ws.Range("L91").Value = 1978.5714
$ws.Range("M91").Value = -903.1999999999998
$ws.Range("N91").Value = -4786.5714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 800
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 13
$ws.Range("N5").Value = -1726

$ws.Range("H86").Value = 4678.7354
$ws.Range("I86").Value = 4005.5454
$ws.Range("J86").Value = 5912.9165
$ws.Range("K86").Value = 4005.5454
$ws.Range("L86").Value = 5912.9165
$ws.Range("M86").Value = -2882.5454
$ws.Range("N86").Value = -8158.9165

$ws.Range("H89").Value = 4678.7354
$ws.Range("I89").Value = 4005.5454
$ws.Range("J89").Value = 5912.9165
$ws.Range("K89").Value = 20027.727
$ws.Range("L89").Value = 29564.5825
$ws.Range("M89").Value = -14411.727
$ws.Range("N89").Value = -40796.5825

$ws.Range("H105").Value = 2244.898
$ws.Range("I105").Value = 2102.5642
$ws.Range("K105").Value = 2102.5642
$ws.Range("M105").Value = -355.5641999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 151002000
$ws.Range("J2").Value = 151002000
$ws.Range("L2").Value = 151002000
$ws.Range("N2").Value = -151002226

$ws.Range("H3").Value = 1126.5
$ws.Range("I3").Value = 250
$ws.Range("J3").Value = 2003
$ws.Range("K3").Value = 250
$ws.Range("L3").Value = 2003
$ws.Range("M3").Value = -137
$ws.Range("N3").Value = -2229

$ws.Range("H4").Value = 205333.33
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 340888.88
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 340888.88
$ws.Range("M4").Value = -1888
$ws.Range("N4").Value = -341112.88

$ws.Range("H5").Value = 187.7
$ws.Range("I5").Value = 99.57143000000001
$ws.Range("J5").Value = 393.33334
$ws.Range("K5").Value = 99.57143000000001
$ws.Range("L5").Value = 393.33334
$ws.Range("M5").Value = 12.42856999999999
$ws.Range("N5").Value = -617.33334

$ws.Range("H8").Value = 2000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2000
$ws.Range("N8").Value = -2280
$ws.Range("M8").ClearContents()

$ws.Range("H10").Value = 312.2
$ws.Range("I10").Value = 312.2
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 312.2
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -173.2
$ws.Range("N10").ClearContents()

$ws.Range("H13").Value = 400
$ws.Range("J13").Value = 400
$ws.Range("L13").Value = 400
$ws.Range("N13").Value = -678

$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

$ws.Range("H15").Value = 571.4286
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 571.4286
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 571.4286
$ws.Range("N15").Value = -911.4286
$ws.Range("M15").ClearContents()

$ws.Range("H16").Value = 832.5
$ws.Range("I16").Value = 777
$ws.Range("J16").Value = 888
$ws.Range("K16").Value = 777
$ws.Range("L16").Value = 888
$ws.Range("M16").Value = -490
$ws.Range("N16").Value = -1462

$ws.Range("H21").Value = 7735.1665
$ws.Range("J21").Value = 7735.1665
$ws.Range("L21").Value = 7735.1665
$ws.Range("N21").Value = -8205.166499999999

$ws.Range("H25").Value = 5802.2
$ws.Range("I25").Value = 3003.6667
$ws.Range("K25").Value = 3003.6667
$ws.Range("M25").Value = -2829.6667

$ws.Range("H107").Value = 883.46155
$ws.Range("I107").Value = 314.16666
$ws.Range("J107").Value = 1371.4286
$ws.Range("K107").Value = 314.16666
$ws.Range("L107").Value = 1371.4286
$ws.Range("M107").Value = 1605.83334
$ws.Range("N107").Value = -5211.4286

$ws.Range("H113").Value = 832.5
$ws.Range("I113").Value = 777
$ws.Range("J113").Value = 888
$ws.Range("K113").Value = 777
$ws.Range("L113").Value = 888
$ws.Range("M113").Value = 1393
$ws.Range("N113").Value = -5228

$ws.Range("H123").Value = 39198
$ws.Range("J123").Value = 39198
$ws.Range("L123").Value = 39198
$ws.Range("N123").Value = -48998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 2757.25
$ws.Range("I119").Value = 514.5
$ws.Range("J119").Value = 5000
$ws.Range("K119").Value = 1543.5
$ws.Range("L119").Value = 15000
$ws.Range("M119").Value = 3294.5
$ws.Range("N119").Value = -24676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4262.2383
$ws.Range("I70").Value = 3903
$ws.Range("J70").Value = 4531.6665
$ws.Range("K70").Value = 3903
$ws.Range("L70").Value = 4531.6665
$ws.Range("M70").Value = -3633
$ws.Range("N70").Value = -5071.6665

$ws.Range("H73").Value = 4262.2383
$ws.Range("I73").Value = 3903
$ws.Range("J73").Value = 4531.6665
$ws.Range("K73").Value = 3903
$ws.Range("L73").Value = 4531.6665
$ws.Range("M73").Value = -2967
$ws.Range("N73").Value = -6403.6665

$ws.Range("H97").Value = 1780.45
$ws.Range("I97").Value = 1104.3636
$ws.Range("J97").Value = 2606.7778
$ws.Range("K97").Value = 1104.3636
$ws.Range("L97").Value = 2606.7778
$ws.Range("M97").Value = -608.3635999999999
$ws.Range("N97").Value = -3598.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 534
$ws.Range("I9").Value = 760
$ws.Range("J9").Value = 383.33334
$ws.Range("K9").Value = 760
$ws.Range("L9").Value = 383.33334
$ws.Range("M9").Value = -536
$ws.Range("N9").Value = -831.33334

$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H12").Value = 991.5
$ws.Range("I12").Value = 1003
$ws.Range("J12").Value = 980
$ws.Range("K12").Value = 1003
$ws.Range("L12").Value = 980
$ws.Range("M12").Value = -833
$ws.Range("N12").Value = -1320

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
